# chore: adapt column header formatting to respective input file names
#
# Renames the paired "_old" / "_new" column header suffixes (columns A-J and
# L-U of the header row) to the concrete format-version identifiers
# "_FV2410" (old/before) and "_FV2504" (new/after), wraps the data range in
# an Excel Table ("Table1"), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# Columns A:J carried the "_old" suffix -> "_FV2410"
# Column  K   is the literal "diff" header (unchanged)
# Columns L:U carried the "_new" suffix -> "_FV2504"

$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldCols) {
    $cell = $ws.Range("$col`1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2410")
}

foreach ($col in $newCols) {
    $cell = $ws.Range("$col`1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2504")
}

# --- 2. Turn the data range into an Excel Table ----------------------------

$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A1:U57"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------------

$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
